$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B answers -------------------------------------------------
# New answers are written in the exact order the authoring session first
# introduced each unique string, so the regenerated shared-string table
# lines up with the source workbook (some answers repeat verbatim and
# reuse an already-known string, e.g. the enterprise/org-process-assets
# input list reused across several "Estimate/Acquire/Develop ... inputs"
# rows).
$ws.Range("B68").Value = "1. Cost estimates`n2. Basis of estimates`n3. Project documents updates"
$ws.Range("B69").Value = "1. Project management plan`n2. Project documents`n3. Business documents`n4. Agreements`n5. Enterprise environmental factors`n6. Organizational process assets"
$ws.Range("B70").Value = "1. Expert judgment`n2. Cost aggregation`n3. Data analysis`n4. Historical information review`n5. Funding limit reconciliation`n6. Financing"
$ws.Range("B71").Value = "1. Cost baseline`n2. Project funding requirements`n3. Project documents updates"
$ws.Range("B72").Value = "1. Project management plan`n2. Project documents`n3. Project funding requirements`n4. Work performance data`n5. Organizational process assets"
$ws.Range("B73").Value = "1. Expert judgment`n2. Data analysis`n3. To-complete performance index`n4. Project management information system"
$ws.Range("B74").Value = "1. Work performance information`n2. Cost forecasts`n3. Change requests`n4. Project management plan updates`n5. Project documents updates"
$ws.Range("B76").Value = "1. Project charter`n2. Project management plan`n3. Project documents`n4. Enterprise environmental factors`n5. Organizational process assets"
$ws.Range("B86").Value = "1. Project charter`n2. Project management plan`n3. Project documents`n4. Enterprise environmental factors`n5. Organizational process assets"
$ws.Range("B105").Value = "1. Project charter`n2. Project management plan`n3. Project documents`n4. Enterprise environmental factors`n5. Organizational process assets"
$ws.Range("B77").Value = "1. Expert judgment`n2. Data gathering`n3. Data analysis`n4. Decision making`n5. Data representation`n6. Test and inspection planning`n7. Meetings"
$ws.Range("B78").Value = "1. Quality management plan`n2. Quality metrics`n3. Project management plan updates`n4. Project documents updates"
$ws.Range("B79").Value = "1. Project management plan`n2. Project documents`n3. Organizational process assets"
$ws.Range("B80").Value = "1. Data gathering`n2. Data analysis`n3. Decision making`n4. Data representation`n5. Audits`n6. Design for X`n7. Problem solving`n8. Quality improvement methods"
$ws.Range("B81").Value = "1. Quality reports`n2. Test and evaluation documents`n3. Change requests`n4. Project management plan updates`n5. Project documents updates"
$ws.Range("B82").Value = "1. Project management plan`n2. Project documents`n3. Approved change requests`n4. Deliverables`n5. Work performance data`n6. Enterprise environmental factors`n7. Organizational process assets"
$ws.Range("B83").Value = "1. Data gathering`n2. Data analysis`n3. Inspection`n4. Testing/product evaluations`n5. Data representation`n6. Meetings"
$ws.Range("B84").Value = "1. Quality control measurements`n2. Verified deliverables`n3. Work performance information`n4. Change requests`n5. Project management plan updates`n6. Project documents updates"
$ws.Range("B87").Value = "1. Expert judgment`n2. Data representation`n3. Organizational theory`n4. Meetings"
$ws.Range("B88").Value = "1. Resources management plan`n2. Team charter`n3. Project documents updates"
$ws.Range("B90").Value = "1. Expert judgment`n2. Bottom-up estimating`n3. Analogous estimating`n4. Parametric estimating`n5. Data analysis`n6. Project management information system`n7. Meetings"
$ws.Range("B91").Value = "1. Resource requirements`n2. Basis of estimates`n3. Resource breakdown structure`n4. Project documents updates"
$ws.Range("B93").Value = "1. Decision making`n2. Interpersonal and team skills`n3. Pre-assignment`n4. Virtual teams"
$ws.Range("B94").Value = "1. Physical resource assignments`n2. Project team assignments`n3. Resource calendars`n4. Change requests`n5. Project management plan updates`n6. Project documents updates`n7. Enterprise environmental factors updates`n8. Organizational process assets updates"
$ws.Range("B96").Value = "1. Colocation`n2. Virtual teams`n3. Communication technology`n4. Interpersonal and team skills`n5. Recognition and rewards`n6. Training`n7. Individual and team assessments`n8. Meetings"
$ws.Range("B97").Value = "1. Team performance assessments`n2. Change requests`n3. Project management plan updates`n4. Project documents updates`n5. Enterprise environmental factors updates`n6. Organizational process assets updates"
$ws.Range("B101").Value = "1. Project management plan`n2. Project documents`n3. Work performance data`n4. Agreements`n5. Organizational process assets"
$ws.Range("B102").Value = "1. Data analysis`n2. Problem solving`n3. Interpersonal and team skills`n4. Project management information system"
$ws.Range("B103").Value = "1. Work performance information`n2. Change requests`n3. Project management plan updates`n4. Project documents updates"
$ws.Range("B113").Value = "1. Work performance information`n2. Change requests`n3. Project management plan updates`n4. Project documents updates"
$ws.Range("B98").Value = "1. Project management plan`n2. Project documents`n3. Work performance reports`n4. Team performance assessments`n5. Enterprise environmental factors`n6. Organizational process assets"
$ws.Range("B99").Value = "1. Interpersonal and team skills`n2. Project management information system"
$ws.Range("B100").Value = "1. Change requests`n2. Project management plan updates`n3. Project documents updates`n4. Enterprise environmental factors updates"
$ws.Range("B106").Value = "1. Expert judgment`n2. Communication requirements analysis`n3. Communication technology`n4. Communication models`n5. Communication methods`n6. Interpersonal and team skills`n7. Data representation`n8. Meetings"
$ws.Range("B107").Value = "1. Communications management plan`n2. Project management plan updates`n3. Project documents update"
$ws.Range("B108").Value = "1. Project management plan`n2. Project documents`n3. Work performance reports`n4. Enterprise environmental factors`n5. Organizational process assets"
$ws.Range("B109").Value = "1. Communication technology`n2. Communication methods`n3. Communication skills`n4. Project management information system`n5. Project reporting`n6. Interpersonal and team skills`n7. Meetings"
$ws.Range("B110").Value = "1. Project communications`n2. Project management plan updates`n3. Project documents updates`n4. Organizational process assets updates"
$ws.Range("B111").Value = "1. Project management plan`n2. Project documents`n3. Work performance data`n4. Enterprise environmental factors`n5. Organizational process assets"
$ws.Range("B112").Value = "1. Expert judgment`n2. Project management information system`n3. Data representation`n4. Interpersonal and team skills`n5. Meetings"

# Rows whose answer text duplicates an answer already present elsewhere
# in the workbook before this edit.
$ws.Range("B89").Value = "1. Project management plan`n2. Project documents`n3. Enterprise environmental factors`n4. Organizational process assets"
$ws.Range("B92").Value = "1. Project management plan`n2. Project documents`n3. Enterprise environmental factors`n4. Organizational process assets"
$ws.Range("B95").Value = "1. Project management plan`n2. Project documents`n3. Enterprise environmental factors`n4. Organizational process assets"

# --- Row heights (auto-sized by Excel for the wrapped answer text) ----
$ws.Rows.Item(68).RowHeight = 45
$ws.Rows.Item(69).RowHeight = 90
$ws.Rows.Item(70).RowHeight = 90
$ws.Rows.Item(71).RowHeight = 45
$ws.Rows.Item(72).RowHeight = 75
$ws.Rows.Item(73).RowHeight = 60
$ws.Rows.Item(74).RowHeight = 75
$ws.Rows.Item(76).RowHeight = 75
$ws.Rows.Item(77).RowHeight = 105
$ws.Rows.Item(78).RowHeight = 60
$ws.Rows.Item(79).RowHeight = 45
$ws.Rows.Item(80).RowHeight = 120
$ws.Rows.Item(81).RowHeight = 75
$ws.Rows.Item(82).RowHeight = 105
$ws.Rows.Item(83).RowHeight = 90
$ws.Rows.Item(84).RowHeight = 90
$ws.Rows.Item(86).RowHeight = 75
$ws.Rows.Item(87).RowHeight = 60
$ws.Rows.Item(88).RowHeight = 45
$ws.Rows.Item(89).RowHeight = 60
$ws.Rows.Item(90).RowHeight = 105
$ws.Rows.Item(91).RowHeight = 60
$ws.Rows.Item(92).RowHeight = 60
$ws.Rows.Item(93).RowHeight = 60
$ws.Rows.Item(94).RowHeight = 120
$ws.Rows.Item(95).RowHeight = 60
$ws.Rows.Item(96).RowHeight = 120
$ws.Rows.Item(97).RowHeight = 90
$ws.Rows.Item(98).RowHeight = 90
$ws.Rows.Item(99).RowHeight = 30
$ws.Rows.Item(100).RowHeight = 60
$ws.Rows.Item(101).RowHeight = 75
$ws.Rows.Item(102).RowHeight = 60
$ws.Rows.Item(103).RowHeight = 60
$ws.Rows.Item(105).RowHeight = 75
$ws.Rows.Item(106).RowHeight = 120
$ws.Rows.Item(107).RowHeight = 45
$ws.Rows.Item(108).RowHeight = 75
$ws.Rows.Item(109).RowHeight = 105
$ws.Rows.Item(110).RowHeight = 60
$ws.Rows.Item(111).RowHeight = 75
$ws.Rows.Item(112).RowHeight = 75
$ws.Rows.Item(113).RowHeight = 60

# --- Restore the active selection to match the authored view state ----
$ws.Range("B115").Select() | Out-Null
